$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("V9").Value = 0.80264970000000002
$ws.Range("W9").Value = 0.8785425
$ws.Range("AB9").Value = 0.80692770000000003
$ws.Range("AT9").Value = 0.84352769999999999
$ws.Range("BN9").Value = 0.86262629999999996
$ws.Range("BS9").Value = 0.87546310000000005
$ws.Range("CM9").Value = 0.87393399999999999
$ws.Range("CX9").Value = 0.888378

$ws.Range("BN10").Value = 0.85194550000000002

$ws.Range("AB11").Value = 0.56446209999999997
$ws.Range("AT11").Value = 0.75286410000000004
$ws.Range("BN11").Value = 0.35143619999999998

$ws.Range("V13").Value = 0.62599470000000002
$ws.Range("W13").Value = 0.83265860000000003
$ws.Range("AB13").Value = 0.691469
$ws.Range("AT13").Value = 0.6728307
$ws.Range("BN13").Value = 0.69969700000000001
$ws.Range("BS13").Value = 0.74074070000000003
$ws.Range("CM13").Value = 0.69070209999999999
$ws.Range("CX13").Value = 0.81660949999999999

